# Updated cryptos list - refresh Price (D) and Volume(1h) (E) columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "'26.497.39"
$ws.Cells.Item(2, 4).Style = "Normal"
$ws.Cells.Item(2, 5).Value = '  +1.44%  '
$ws.Cells.Item(3, 4).Value = "'1.678.28"
$ws.Cells.Item(3, 4).Style = "Normal"
$ws.Cells.Item(3, 5).Value = '  +1.86%  '
$ws.Cells.Item(4, 4).Value = "'1.002"
$ws.Cells.Item(4, 4).Style = "Normal"
$ws.Cells.Item(4, 5).Value = '  +0.05%  '
$ws.Cells.Item(5, 4).Value = "'219.45"
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).Value = '  +1.90%  '
$ws.Cells.Item(6, 4).Value = "'0.5311"
$ws.Cells.Item(6, 4).Style = "Normal"
$ws.Cells.Item(6, 5).Value = '  +1.48%  '
$ws.Cells.Item(7, 4).Value = "'1.002"
$ws.Cells.Item(7, 4).Style = "Normal"
$ws.Cells.Item(7, 5).Value = '  +0.10%  '
$ws.Cells.Item(8, 4).Value = "'0.2700"
$ws.Cells.Item(8, 4).Style = "Normal"
$ws.Cells.Item(8, 5).Value = '  +3.39%  '
$ws.Cells.Item(9, 4).Value = "'0.06411"
$ws.Cells.Item(9, 4).Style = "Normal"
$ws.Cells.Item(9, 5).Value = '  +1.37%  '
$ws.Cells.Item(10, 4).Value = "'21.83"
$ws.Cells.Item(10, 4).Style = "Normal"
$ws.Cells.Item(10, 5).Value = '  +4.79%  '
$ws.Cells.Item(11, 4).Value = "'0.07806"
$ws.Cells.Item(11, 4).Style = "Normal"
$ws.Cells.Item(11, 5).Value = '  +1.46%  '
$ws.Cells.Item(12, 4).Value = "'1.689.68"
$ws.Cells.Item(12, 4).Style = "Normal"
$ws.Cells.Item(12, 5).Value = '  +2.49%  '
$ws.Cells.Item(13, 4).Value = "'4.514"
$ws.Cells.Item(13, 4).Style = "Normal"
$ws.Cells.Item(13, 5).Value = '  +2.01%  '
$ws.Cells.Item(14, 4).Value = "'0.5595"
$ws.Cells.Item(14, 4).Style = "Normal"
$ws.Cells.Item(14, 5).Value = '  +0.20%  '
$ws.Cells.Item(15, 4).Value = "'0.0₅8352"
$ws.Cells.Item(15, 4).Style = "Normal"
$ws.Cells.Item(15, 5).Value = '  +1.77%  '
$ws.Cells.Item(16, 4).Value = "'65.76"
$ws.Cells.Item(16, 4).Style = "Normal"
$ws.Cells.Item(16, 5).Value = '  +0.81%  '
$ws.Cells.Item(17, 4).Value = "'26.538.85"
$ws.Cells.Item(17, 4).Style = "Normal"
$ws.Cells.Item(17, 5).Value = '  +1.49%  '
$ws.Cells.Item(18, 5).Value = '  +0.01%  '
$ws.Cells.Item(19, 4).Value = "'4.800"
$ws.Cells.Item(19, 4).Style = "Normal"
$ws.Cells.Item(19, 5).Value = '  +1.07%  '
$ws.Cells.Item(20, 4).Value = "'193.36"
$ws.Cells.Item(20, 4).Style = "Normal"
$ws.Cells.Item(20, 5).Value = '  +2.14%  '
$ws.Cells.Item(21, 4).Value = "'10.33"
$ws.Cells.Item(21, 4).Style = "Normal"
$ws.Cells.Item(21, 5).Value = '  +0.80%  '
$ws.Cells.Item(22, 4).Value = "'6.343"
$ws.Cells.Item(22, 4).Style = "Normal"
$ws.Cells.Item(22, 5).Value = '  +2.18%  '
$ws.Cells.Item(23, 5).Value = '  +0.14%  '
$ws.Cells.Item(24, 4).Value = "'0.1275"
$ws.Cells.Item(24, 4).Style = "Normal"
$ws.Cells.Item(24, 5).Value = '  +5.37%  '
$ws.Cells.Item(25, 4).Value = "'138.76"
$ws.Cells.Item(25, 4).Style = "Normal"
$ws.Cells.Item(25, 5).Value = '  -4.77%  '
$ws.Cells.Item(26, 4).Value = "'7.417"
$ws.Cells.Item(26, 4).Style = "Normal"
$ws.Cells.Item(26, 5).Value = '  -0.31%  '
$ws.Cells.Item(27, 5).Value = '  +2.32%  '
$ws.Cells.Item(28, 4).Value = "'1.443"
$ws.Cells.Item(28, 4).Style = "Normal"
$ws.Cells.Item(28, 5).Value = '  +3.50%  '
$ws.Cells.Item(29, 4).Value = "'0.06284"
$ws.Cells.Item(29, 4).Style = "Normal"
$ws.Cells.Item(29, 5).Value = '  +6.66%  '
$ws.Cells.Item(30, 4).Value = "'1.289"
$ws.Cells.Item(30, 4).Style = "Normal"
$ws.Cells.Item(30, 5).Value = '  +2.16%  '
$ws.Cells.Item(31, 4).Value = "'3.608"
$ws.Cells.Item(31, 4).Style = "Normal"
$ws.Cells.Item(31, 5).Value = '  +4.69%  '
$ws.Cells.Item(32, 4).Value = "'3.410"
$ws.Cells.Item(32, 4).Style = "Normal"
$ws.Cells.Item(32, 5).Value = '  -0.07%  '
$ws.Cells.Item(33, 4).Value = "'1.698"
$ws.Cells.Item(33, 4).Style = "Normal"
$ws.Cells.Item(33, 5).Value = '  +2.49%  '
$ws.Cells.Item(34, 4).Value = "'1.014"
$ws.Cells.Item(34, 4).Style = "Normal"
$ws.Cells.Item(34, 5).Value = '  +2.73%  '
$ws.Cells.Item(35, 4).Value = "'0.6194"
$ws.Cells.Item(35, 4).Style = "Normal"
$ws.Cells.Item(35, 5).Value = '  +8.98%  '
$ws.Cells.Item(36, 4).Value = "'2.423"
$ws.Cells.Item(36, 4).Style = "Normal"
$ws.Cells.Item(36, 5).Value = '  +1.38%  '
$ws.Cells.Item(37, 4).Value = "'2.791"
$ws.Cells.Item(37, 4).Style = "Normal"
$ws.Cells.Item(37, 5).Value = '  +0.85%  '
$ws.Cells.Item(38, 4).Value = "'0.01634"
$ws.Cells.Item(38, 4).Style = "Normal"
$ws.Cells.Item(38, 5).Value = '  +0.82%  '
$ws.Cells.Item(39, 4).Value = "'6.124"
$ws.Cells.Item(39, 4).Style = "Normal"
$ws.Cells.Item(39, 5).Value = '  +6.29%  '
$ws.Cells.Item(40, 4).Value = "'1.096.58"
$ws.Cells.Item(40, 4).Style = "Normal"
$ws.Cells.Item(40, 5).Value = '  +6.38%  '
$ws.Cells.Item(41, 4).Value = "'0.8639"
$ws.Cells.Item(41, 4).Style = "Normal"
$ws.Cells.Item(41, 5).Value = '  +0.67%  '
$ws.Cells.Item(42, 4).Value = "'1.000"
$ws.Cells.Item(42, 4).Style = "Normal"
$ws.Cells.Item(42, 5).Value = '  -0.03%  '
$ws.Cells.Item(43, 4).Value = "'100.61"
$ws.Cells.Item(43, 4).Style = "Normal"
$ws.Cells.Item(43, 5).Value = '  +0.39%  '
$ws.Cells.Item(44, 4).Value = "'1.824.07"
$ws.Cells.Item(44, 4).Style = "Normal"
$ws.Cells.Item(44, 5).Value = '  +1.51%  '
$ws.Cells.Item(45, 4).Value = "'0.0₈112"
$ws.Cells.Item(45, 4).Style = "Normal"
$ws.Cells.Item(45, 5).Value = '  +6.78%  '
$ws.Cells.Item(46, 4).Value = "'58.78"
$ws.Cells.Item(46, 4).Style = "Normal"
$ws.Cells.Item(47, 4).Value = "'8.171"
$ws.Cells.Item(47, 4).Style = "Normal"
$ws.Cells.Item(47, 5).Value = '  +0.50%  '
$ws.Cells.Item(48, 4).Value = "'1.003"
$ws.Cells.Item(48, 4).Style = "Normal"
$ws.Cells.Item(48, 5).Value = '  -0.27%  '
$ws.Cells.Item(49, 4).Value = "'1.493"
$ws.Cells.Item(49, 4).Style = "Normal"
$ws.Cells.Item(49, 5).Value = '  +7.27%  '
$ws.Cells.Item(50, 4).Value = "'0.05198"
$ws.Cells.Item(50, 4).Style = "Normal"
$ws.Cells.Item(50, 5).Value = '  +0.72%  '
$ws.Cells.Item(51, 4).Value = "'6.037"
$ws.Cells.Item(51, 4).Style = "Normal"
$ws.Cells.Item(51, 5).Value = '  +1.58%  '
